# Update column G ("K") values on the active sheet to reflect the
# regenerated strike-count (s_vals) data described in the commit message:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals". Only column G (rows 2-49) values change; row 37, 43 and
# 47 already hold 0 and remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 0
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 4
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 2
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    29 = 1
    30 = 0
    31 = 2
    32 = 1
    33 = 2
    34 = 2
    35 = 0
    36 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 1
    42 = 1
    44 = 1
    45 = 0
    46 = 2
    48 = 1
    49 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
